# Apply crypto price/volume updates to columns D (Price) and E (Volume 1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.986.87"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "3.570.02"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.22"
$ws.Range("E5").Value = "  -3.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.93"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("D7").Value = "3.564.31"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("E8").Value = "  -4.65%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.649"
$ws.Range("E11").Value = "  -4.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.00"
$ws.Range("E12").Value = "  -5.65%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.51"
$ws.Range("E14").Value = "  -4.90%  "
$ws.Range("D15").Value = "4.137.68"
$ws.Range("E15").Value = "  -2.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.59"
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("D17").Value = "3.564.23"
$ws.Range("E17").Value = "  -2.53%  "
$ws.Range("D18").Value = "69.858.02"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.56"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("E21").Value = "  -3.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "499.04"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.99"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E24").Value = "  -7.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.37"
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.12"
$ws.Range("E26").Value = "  +3.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.27"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.95"
$ws.Range("E28").Value = "  -7.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.29"
$ws.Range("E29").Value = "  -3.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.60"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.56"
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "66.90"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.04"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("E34").Value = "  -6.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "567.05"
$ws.Range("E35").Value = "  -10.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.13"
$ws.Range("E36").Value = "  +9.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.59"
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "0.0₃0790"
$ws.Range("E39").Value = "  -5.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.394"
$ws.Range("E40").Value = "  -4.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("E43").Value = "  -10.09%  "
$ws.Range("E44").Value = "  -5.54%  "
$ws.Range("D45").Value = "3.223.77"
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.47"
$ws.Range("E46").Value = "  +4.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0440"
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  -3.45%  "
